{"js": "// The document contains a single-column table where each row holds one\n// benchmark statistic. This edit refreshes the stats: some single values\n// change, a handful of new one-value rows appear in place of what used to\n// be a single combined multi-tab summary line, and the three trailing\n// \"combined\" rows (which packed 10 tab-separated numbers into one run)\n// collapse down to just their leading summary number.\n//\n// Net effect: the table keeps the same row count (46), only the text\n// inside certain cells changes. We read the whole table as a 2D array,\n// patch the rows that differ, and write the whole array back in one shot.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\n// New text, by row index (0-based), for every row whose content changes.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"314\",\n  4: \"0.00002\",\n  5: \"0.00016\",\n  6: \"0.00005\",\n  7: \"0.00002\",\n  8: \"0.00004\",\n  9: \"0.00004\",\n  10: \"0.00009\",\n  11: \"0.01558\",\n  43: \"99.98\",\n  44: \"0.02\",\n  45: \"83\",\n};\n\nconst newValues = table.values.map((row) => row.slice());\nfor (const [idxStr, text] of Object.entries(updates)) {\n  const idx = Number(idxStr);\n  newValues[idx] = [text];\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document body is a single-column table where each row holds one\n# benchmark statistic. This edit refreshes the stats: a few single values\n# change, a handful of new one-value rows appear where there used to be a\n# single combined multi-tab summary line, and the three trailing \"combined\"\n# rows (which packed 10 tab-separated numbers into one run) collapse down\n# to just their leading summary number.\n#\n# Net effect: the table keeps the same row count (46); only the text inside\n# certain cells changes. We set each changed cell's Range.Text directly,\n# which also collapses any pre-existing multi-run/tab content in that cell\n# down to the single new run.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"314\"\n    5  = \"0.00002\"\n    6  = \"0.00016\"\n    7  = \"0.00005\"\n    8  = \"0.00002\"\n    9  = \"0.00004\"\n    10 = \"0.00004\"\n    11 = \"0.00009\"\n    12 = \"0.01558\"\n    44 = \"99.98\"\n    45 = \"0.02\"\n    46 = \"83\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
